# Add one new test case on Page Product assert:
# rename Sheet3 -> 004_AccessoriesPage, populate it with the accessories
# product/price table, and adjust workbook/sheet view state to match.

$wb = $excel.ActiveWorkbook

# --- Rename the third sheet and fill in its data -------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Name = "004_AccessoriesPage"

$ws3.Range("A1").Value = "TestCase"
$ws3.Range("A2").Value = "004-Test Accessories Page"

$ws3.Range("B1").Value = "Product Name"
$ws3.Range("B2").Value = "Magic Mouse"
$ws3.Range("B3").Value = "Apple TV"
$ws3.Range("B4").Value = "Sennheiser RS 120"
$ws3.Range("B5").Value = "Skullcandy PLYR 1 – Black"
$ws3.Range("B6").Value = "Apple 27 inch Thunderbolt Display"
$ws3.Range("B7").Value = "Asus MX239H 23-inch Widescreen AH"

$ws3.Range("C1").Value = "Product Prices"

$ws3.Range("C2:C7").NumberFormat = "@"
$ws3.Range("C2").Value = "$150.00"
$ws3.Range("C3").Value = "$80.00"
$ws3.Range("C4").Value = "$50.00"
$ws3.Range("C5").Value = "$110.00"
$ws3.Range("C6").Value = "$764.00"
$ws3.Range("C7").Value = "$199.00"

$ws3.Columns.Item(1).ColumnWidth = 22.6640625
$ws3.Columns.Item(2).ColumnWidth = 31.77734375
$ws3.Columns.Item(3).ColumnWidth = 12.5546875

# Match sheet2's print setup (same paper/orientation) on the new sheet.
$ps3 = $ws3.PageSetup
$ps3.PaperSize = 9
$ps3.Orientation = 1

# --- Sheet view / selection bookkeeping -----------------------------------
$ws1 = $wb.Worksheets.Item("001_LoginCorrectly")
$ws1.Range("A14").Select() | Out-Null

$ws2 = $wb.Worksheets.Item("002_LoginIncorrectly")
$ws2.Range("A1:C3").Select() | Out-Null

$ws3.Range("C6").Select() | Out-Null

# Make the accessories sheet the active tab, scrolled so it's the first
# visible sheet on the tab strip.
$ws3.Activate() | Out-Null
$wb.Windows.Item(1).ScrollWorkbookTabs(1) | Out-Null
